# v2.6 Added decoupled suspension, four-wheel steering, scripts to generate GGV diagram
#
# Adds a new worksheet "Semi_Truck_Scalable" (a clone of "Truck_Amandla_3Axle"),
# makes it the active sheet/tab, updates its "Comments" cell (H3) to reference
# its own name, and nudges the remembered cell selections on the other sheets
# the way Excel does when the active sheet/cell changes.

$wb = $excel.ActiveWorkbook

# --- Clone "Truck_Amandla_3Axle" to create the new "Semi_Truck_Scalable" sheet ---
$truck = $wb.Worksheets.Item("Truck_Amandla_3Axle")
[void]$truck.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Semi_Truck_Scalable"

# The "Comments" column (H3) on every sheet names the sheet itself.
$newSheet.Range("H3").Value = "Semi_Truck_Scalable"

# --- Truck_Amandla_3Axle keeps its last remembered selection, now at D24 ---
[void]$truck.Activate()
[void]$truck.Range("D24").Select()

# --- The new sheet becomes the active tab, selection parked at J17 ---
[void]$newSheet.Activate()
[void]$newSheet.Range("J17").Select()
